$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.962.80'
$ws.Range('E2').Value = '  -3.57%  '
$ws.Range('D3').Value = '3.836.34'
$ws.Range('E3').Value = '  -3.05%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.21'
$ws.Range('E5').Value = '  -1.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.86'
$ws.Range('E6').Value = '  -2.90%  '
$ws.Range('D7').Value = '3.836.03'
$ws.Range('E7').Value = '  -3.04%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.528'
$ws.Range('E9').Value = '  -1.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.163'
$ws.Range('E10').Value = '  -4.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.42'
$ws.Range('E11').Value = '  -0.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.456'
$ws.Range('E12').Value = '  -3.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000261'
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.84'
$ws.Range('E14').Value = '  -5.19%  '
$ws.Range('D15').Value = '4.492.85'
$ws.Range('E15').Value = '  -2.70%  '
$ws.Range('D16').Value = '3.846.24'
$ws.Range('E16').Value = '  -3.03%  '
$ws.Range('D17').Value = '68.220.86'
$ws.Range('E17').Value = '  -2.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.06'
$ws.Range('E18').Value = '  -1.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.36'
$ws.Range('E19').Value = '  -4.15%  '
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.90'
$ws.Range('E21').Value = '  -1.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '465.58'
$ws.Range('E22').Value = '  -6.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.730'
$ws.Range('E23').Value = '  -2.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000159'
$ws.Range('E24').Value = '  -5.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.82'
$ws.Range('E25').Value = '  -3.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.23'
$ws.Range('E26').Value = '  -4.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.03'
$ws.Range('E27').Value = '  -4.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.00'
$ws.Range('E28').Value = '  -3.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.96'
$ws.Range('E30').Value = '  -2.47%  '
$ws.Range('D31').Value = '3.996.44'
$ws.Range('E31').Value = '  -2.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.56'
$ws.Range('E32').Value = '  -4.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.30'
$ws.Range('E33').Value = '  -6.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.05'
$ws.Range('E34').Value = '  -4.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.56'
$ws.Range('E35').Value = '  -1.49%  '
$ws.Range('D36').Value = '3.809.31'
$ws.Range('E36').Value = '  -2.81%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.104'
$ws.Range('E37').Value = '  -4.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.57'
$ws.Range('E38').Value = '  +7.53%  '
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.88'
$ws.Range('E41').Value = '  -4.82%  '
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.312'
$ws.Range('E43').Value = '  -5.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.000296'
$ws.Range('E44').Value = '  +4.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.97'
$ws.Range('E45').Value = '  -7.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '419.83'
$ws.Range('E46').Value = '  -4.98%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.59'
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '46.94'
$ws.Range('E49').Value = '  -2.97%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '141.81'
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.10'
$ws.Range('E51').Value = '  +1.34%  '
